$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for worker CARMELITA PAJARO POSSO (45469039) entirely -
# this period/worker left the account statement in this update.
$ws.Range("B16").EntireRow.Delete()

# Rewrite the remaining detail rows (now shifted up to rows 16-35) with the
# refreshed / re-sorted data for "parte 1 de nuevos estado de cuenta".
$data = @(
  @(16, "CC", "1026256860", "TANIA DEL CARMEN PUENTE ALARCON", "2102", 161120, 4028000),
  @(17, "CC", "1026256860", "TANIA DEL CARMEN PUENTE ALARCON", "2103", 161120, 4028000),
  @(18, "CC", "1026256860", "TANIA DEL CARMEN PUENTE ALARCON", "2104", 161120, 4028000),
  @(19, "CC", "1026256860", "TANIA DEL CARMEN PUENTE ALARCON", "2105", 161120, 4028000),
  @(20, "CC", "1143374150", "GUSTAVO ADOLFO MONSALVE JIMENEZ", "2309", 6149, 4611657),
  @(21, "CC", "1000913551", "ALEJANDRA HERNANDEZ AGUIRRE", "2404", 1851, 1388152),
  @(22, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2407", 20800, 1300000),
  @(23, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2408", 52000, 1300000),
  @(24, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2409", 52000, 1300000),
  @(25, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2410", 52000, 1300000),
  @(26, "CC", "1235041215", "LAURA VANESSA HERNANDEZ BAENA", "2411", 12298, 4611657),
  @(27, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2411", 52000, 1300000),
  @(28, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2412", 52000, 1300000),
  @(29, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2501", 52000, 1300000),
  @(30, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2502", 52000, 1300000),
  @(31, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2503", 52000, 1300000),
  @(32, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2504", 52000, 1300000),
  @(33, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2505", 52000, 1300000),
  @(34, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2506", 52000, 1300000),
  @(35, "CC", "1143367691", "MAIRA ALEJANDRA PUELLO DE HOYOS", "2507", 52000, 1300000)
)

foreach ($row in $data) {
  $r = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $ws.Cells.Item($r, 4).Value = $row[3]
  $ws.Cells.Item($r, 5).Value = $row[4]
  $ws.Cells.Item($r, 6).Value = $row[5]
  $ws.Cells.Item($r, 7).Value = $row[6]
}

# Refresh the summary totals above the table: total overdue value and
# worker headcount (period count stays the same at 19).
$ws.Range("E11").Value = 1309578
$ws.Range("C13").Value = 5
